$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# Add a new worksheet right after the existing "Sheet" and name it.
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheet1)
$newSheet.Name = "New Sheet"

# Enter the "Nome"/"CPF" table. Headers and names are typed first (column A),
# then the CPF numbers are filled in out of row-order (row 4, then 2, then 3).
$newSheet.Range("A1").Value = "Nome"
$newSheet.Range("B1").Value = "CPF"
$newSheet.Range("A2").Value = "Anilton"
$newSheet.Range("A3").Value = "Quixote"
$newSheet.Range("A4").Value = "Laura"
$newSheet.Range("B4").Value = "156.874.960-20"
$newSheet.Range("B2").Value = "818.965.205-25"
$newSheet.Range("B3").Value = "202.350.365-96"

# Widen column B so the CPF values fit.
$newSheet.Columns.Item(2).ColumnWidth = 13.14

# Match Excel's standard page margins on the new sheet.
$newSheet.PageSetup.LeftMargin = $excel.InchesToPoints(0.75)
$newSheet.PageSetup.RightMargin = $excel.InchesToPoints(0.75)
$newSheet.PageSetup.TopMargin = $excel.InchesToPoints(1)
$newSheet.PageSetup.BottomMargin = $excel.InchesToPoints(1)
$newSheet.PageSetup.HeaderMargin = $excel.InchesToPoints(0.5)
$newSheet.PageSetup.FooterMargin = $excel.InchesToPoints(0.5)

# Leave the cursor resting on column C, after the entered data.
$newSheet.Columns("C").Select() | Out-Null
